$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the connection schema: re-number the free/unused Arduino pins
# (previously mislabeled as PIN #9 / #10 / #11 / #12 / #13) to the
# correct PIN #3 / #4 / #5 / #6 / #7.
$ws.Range("C4").Value = "PIN #3"
$ws.Range("C5").Value = "PIN #4"
$ws.Range("C6").Value = "PIN #5"
$ws.Range("C3").Value = "PIN #6"
$ws.Range("C2").Value = "PIN #7"

# Move the active selection to C2
$ws.Range("C2").Select()
